$d = $word.ActiveDocument

# 1. In-text citation year updates: Oyafuso 2024 -> 2025 (both occurrences)
$d.Content.Find.Execute("(Oyafuso, 2024)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(Oyafuso, 2025)", 2)

# 2. In-text citation year updates: Rohan 2024 -> 2025 (both occurrences)
$d.Content.Find.Execute("(Rohan, 2024)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(Rohan, 2025)", 2)

# 3. Figure 1 alt text (docPr descr) - add "(kilotons)"
$d.InlineShapes.Item(1).AlternativeText = "Figure 1. Biomass index (kilotons) of jellyfish from RACE Groundfish Assessment Program summer bottom trawl surveys of the Aleutian Islands from 1991 to 2024 showing the observed survey biomass index mean (blue points), random effects model fitted mean (solid black line), 95% confidence interval (gray shading), overall time series mean (solid gray line), and horizontal dashed gray lines representing one standard deviation from the mean."

# 4. Figure 1 caption text - add "(kilotons)"
$d.Content.Find.Execute("Figure 1. Biomass index of jellyfish from RACE Groundfish Assessment Program summer bottom trawl surveys of the Aleutian Islands from 1991 to 2024 showing the observed survey biomass index mean (blue points), random effects model fitted mean (solid black line), 95% confidence interval (gray shading), overall time series mean (solid gray line), and horizontal dashed gray lines representing one standard deviation from the mean.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Figure 1. Biomass index (kilotons) of jellyfish from RACE Groundfish Assessment Program summer bottom trawl surveys of the Aleutian Islands from 1991 to 2024 showing the observed survey biomass index mean (blue points), random effects model fitted mean (solid black line), 95% confidence interval (gray shading), overall time series mean (solid gray line), and horizontal dashed gray lines representing one standard deviation from the mean.", 2)

# 5. Figure 2 alt text (docPr descr) - add "(kilotons)" and change wording at the end
$d.InlineShapes.Item(2).AlternativeText = "Figure 2. Biomass index (kilotons) of jellyfish in Aleutian Islands subareas (Southern Bering Sea [SBS], Eastern Aleutian Islands [EAI], Central Aleutian Islands [CAI], and Western Aleutian Islands [WAI]) estimated from RACE Groundfish Assessment Program summer bottom trawl survey data from 1991 to 2024."

# 6. Figure 2 caption text - add "(kilotons)" and change wording at the end
$d.Content.Find.Execute("Figure 2. Biomass index of jellyfish in Aleutian Islands subareas (Southern Bering Sea [SBS], Eastern Aleutian Islands [EAI], Central Aleutian Islands [CAI], and Western Aleutian Islands [WAI]) estimated from RACE Groundfish Assessment Program summer bottom trawl surveys of the Aleutian Islands from 1991 to 2024.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Figure 2. Biomass index (kilotons) of jellyfish in Aleutian Islands subareas (Southern Bering Sea [SBS], Eastern Aleutian Islands [EAI], Central Aleutian Islands [CAI], and Western Aleutian Islands [WAI]) estimated from RACE Groundfish Assessment Program summer bottom trawl survey data from 1991 to 2024.", 2)

# 7. Bibliography: Oyafuso year and package version
$d.Content.Find.Execute("Oyafuso, Z. (2024).", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Oyafuso, Z. (2025).", 2)
# NOTE: search text is a substring inside the italic run (not the full run span)
# so the run's formatting (italics) is preserved across the edit.
$d.Content.Find.Execute("R package version 3.0.0", $true, $false, $false, $false, $false,
                         $true, 1, $false, "R package version 3.0.2", 2)

# 8. Bibliography: Rohan year and package version
$d.Content.Find.Execute("Rohan, S. (2024).", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Rohan, S. (2025).", 2)
$d.Content.Find.Execute("R package version 0.1.0", $true, $false, $false, $false, $false,
                         $true, 1, $false, "R package version 1.2.0", 2)
